$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 67952
$ws.Range("E2").Value = 1339820854507
$ws.Range("F2").Value = 12027414010
$ws.Range("G2").Value = 0.44596

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 3791.33
$ws.Range("E3").Value = 455431131112
$ws.Range("F3").Value = 8790199617
$ws.Range("G3").Value = -0.03487

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 0.99915
$ws.Range("E4").Value = 112074706805
$ws.Range("F4").Value = 18644294445
$ws.Range("G4").Value = 0.00257

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 599.73
$ws.Range("E5").Value = 92288356635
$ws.Range("F5").Value = 835855298
$ws.Range("G5").Value = 0.6885

$ws.Range("B6").Value = "SOL"
$ws.Range("C6").Value = "Solana"
$ws.Range("D6").Value = 165.34
$ws.Range("E6").Value = 75966803423
$ws.Range("F6").Value = 1109749884
$ws.Range("G6").Value = -0.86615

$ws.Range("B7").Value = "STETH"
$ws.Range("C7").Value = "Lido Staked Ether"
$ws.Range("D7").Value = 3789.67
$ws.Range("E7").Value = 35993923577
$ws.Range("F7").Value = 38275320
$ws.Range("G7").Value = -0.0379

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "USDC"
$ws.Range("D8").Value = 0.9998860000000001
$ws.Range("E8").Value = 32350141391
$ws.Range("F8").Value = 2433705617
$ws.Range("G8").Value = -0.01199

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "XRP"
$ws.Range("D9").Value = 0.516373
$ws.Range("E9").Value = 28613555361
$ws.Range("F9").Value = 455041709
$ws.Range("G9").Value = -0.67903

$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "Dogecoin"
$ws.Range("D10").Value = 0.158682
$ws.Range("E10").Value = 22899880120
$ws.Range("F10").Value = 499788700
$ws.Range("G10").Value = -0.43382

$ws.Range("B11").Value = "ADA"
$ws.Range("C11").Value = "Cardano"
$ws.Range("D11").Value = 0.450885
$ws.Range("E11").Value = 15962072933
$ws.Range("F11").Value = 183123806
$ws.Range("G11").Value = 0.5378500000000001

$ws.Range("B12").Value = "TON"
$ws.Range("C12").Value = "Toncoin"
$ws.Range("D12").Value = 6.43
$ws.Range("E12").Value = 15489964987
$ws.Range("F12").Value = 221780480
$ws.Range("G12").Value = 0.89371

$ws.Range("B13").Value = "SHIB"
$ws.Range("C13").Value = "Shiba Inu"
$ws.Range("D13").Value = 0.0000249
$ws.Range("E13").Value = 14659811169
$ws.Range("F13").Value = 314260975
$ws.Range("G13").Value = -1.43707

$ws.Range("B14").Value = "AVAX"
$ws.Range("C14").Value = "Avalanche"
$ws.Range("D14").Value = 35.87
$ws.Range("E14").Value = 14117469234
$ws.Range("F14").Value = 181205307
$ws.Range("G14").Value = -0.36571

$ws.Range("B15").Value = "LINK"
$ws.Range("C15").Value = "Chainlink"
$ws.Range("D15").Value = 18.27
$ws.Range("E15").Value = 10730453044
$ws.Range("F15").Value = 282502776
$ws.Range("G15").Value = -1.71004

$ws.Range("B16").Value = "WBTC"
$ws.Range("C16").Value = "Wrapped Bitcoin"
$ws.Range("D16").Value = 67964
$ws.Range("E16").Value = 10597510029
$ws.Range("F16").Value = 105325218
$ws.Range("G16").Value = 0.35564

$ws.Range("B17").Value = "TRX"
$ws.Range("C17").Value = "TRON"
$ws.Range("D17").Value = 0.113931
$ws.Range("E17").Value = 9946431798
$ws.Range("F17").Value = 215966288
$ws.Range("G17").Value = 1.81523

$ws.Range("B18").Value = "DOT"
$ws.Range("C18").Value = "Polkadot"
$ws.Range("D18").Value = 7.07
$ws.Range("E18").Value = 9681551497
$ws.Range("F18").Value = 114845581
$ws.Range("G18").Value = -0.1377

$ws.Range("B19").Value = "BCH"
$ws.Range("C19").Value = "Bitcoin Cash"
$ws.Range("D19").Value = 461.4
$ws.Range("E19").Value = 9103498570
$ws.Range("F19").Value = 142833360
$ws.Range("G19").Value = 0.37673

$ws.Range("B20").Value = "NEAR"
$ws.Range("C20").Value = "NEAR Protocol"
$ws.Range("D20").Value = 7.44
$ws.Range("E20").Value = 8035867171
$ws.Range("F20").Value = 187581016
$ws.Range("G20").Value = 3.02189

$ws.Range("B21").Value = "UNI"
$ws.Range("C21").Value = "Uniswap"
$ws.Range("D21").Value = 9.73
$ws.Range("E21").Value = 7329448130
$ws.Range("F21").Value = 144247452
$ws.Range("G21").Value = -2.11288

$ws.Range("B22").Value = "MATIC"
$ws.Range("C22").Value = "Polygon"
$ws.Range("D22").Value = 0.695935
$ws.Range("E22").Value = 6454468953
$ws.Range("F22").Value = 186914552
$ws.Range("G22").Value = -0.55875

$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "Pepe"
$ws.Range("D23").Value = 0.00001478
$ws.Range("E23").Value = 6208897685
$ws.Range("F23").Value = 835937480
$ws.Range("G23").Value = -3.50396

$ws.Range("B24").Value = "LTC"
$ws.Range("C24").Value = "Litecoin"
$ws.Range("D24").Value = 82.76000000000001
$ws.Range("E24").Value = 6167812739
$ws.Range("F24").Value = 185194200
$ws.Range("G24").Value = -0.57738

$ws.Range("B25").Value = "ICP"
$ws.Range("C25").Value = "Internet Computer"
$ws.Range("D25").Value = 12.04
$ws.Range("E25").Value = 5584233709
$ws.Range("F25").Value = 42712983
$ws.Range("G25").Value = -0.03594

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "LEO Token"
$ws.Range("D26").Value = 5.92
$ws.Range("E26").Value = 5473902884
$ws.Range("F26").Value = 642839
$ws.Range("G26").Value = -0.00467

$ws.Range("B27").Value = "DAI"
$ws.Range("C27").Value = "Dai"
$ws.Range("D27").Value = 0.999184
$ws.Range("E27").Value = 5312543984
$ws.Range("F27").Value = 204981462
$ws.Range("G27").Value = -0.03749

$ws.Range("B28").Value = "FET"
$ws.Range("C28").Value = "Fetch.ai"
$ws.Range("D28").Value = 2.11
$ws.Range("E28").Value = 5312347307
$ws.Range("F28").Value = 97633423
$ws.Range("G28").Value = 0.37264

$ws.Range("B29").Value = "WEETH"
$ws.Range("C29").Value = "Wrapped eETH"
$ws.Range("D29").Value = 3938.98
$ws.Range("E29").Value = 5107325033
$ws.Range("F29").Value = 41455330
$ws.Range("G29").Value = -0.13055

$ws.Range("B30").Value = "ETC"
$ws.Range("C30").Value = "Ethereum Classic"
$ws.Range("D30").Value = 29.24
$ws.Range("E30").Value = 4299362771
$ws.Range("F30").Value = 117782146
$ws.Range("G30").Value = -1.19622

$ws.Range("B31").Value = "APT"
$ws.Range("C31").Value = "Aptos"
$ws.Range("D31").Value = 9
$ws.Range("E31").Value = 3935174330
$ws.Range("F31").Value = 98322229
$ws.Range("G31").Value = -0.53955

$ws.Range("B32").Value = "RNDR"
$ws.Range("C32").Value = "Render"
$ws.Range("D32").Value = 9.949999999999999
$ws.Range("E32").Value = 3861180088
$ws.Range("F32").Value = 78958546
$ws.Range("G32").Value = -0.35211

$ws.Range("B33").Value = "EZETH"
$ws.Range("C33").Value = "Renzo Restaked ETH"
$ws.Range("D33").Value = 3737.97
$ws.Range("E33").Value = 3724119372
$ws.Range("F33").Value = 21078804
$ws.Range("G33").Value = 0.06818

$ws.Range("B34").Value = "HBAR"
$ws.Range("C34").Value = "Hedera"
$ws.Range("D34").Value = 0.09962
$ws.Range("E34").Value = 3561302018
$ws.Range("F34").Value = 39588529
$ws.Range("G34").Value = -0.28233

$ws.Range("B35").Value = "KAS"
$ws.Range("C35").Value = "Kaspa"
$ws.Range("D35").Value = 0.139011
$ws.Range("E35").Value = 3300630217
$ws.Range("F35").Value = 13861324
$ws.Range("G35").Value = 1.04621

$ws.Range("B36").Value = "WIF"
$ws.Range("C36").Value = "dogwifhat"
$ws.Range("D36").Value = 3.27
$ws.Range("E36").Value = 3273736595
$ws.Range("F36").Value = 344604691
$ws.Range("G36").Value = -2.30103

$ws.Range("B37").Value = "IMX"
$ws.Range("C37").Value = "Immutable"
$ws.Range("D37").Value = 2.21
$ws.Range("E37").Value = 3261896696
$ws.Range("F37").Value = 37020514
$ws.Range("G37").Value = -1.58343

$ws.Range("B38").Value = "ATOM"
$ws.Range("C38").Value = "Cosmos Hub"
$ws.Range("D38").Value = 8.359999999999999
$ws.Range("E38").Value = 3256192618
$ws.Range("F38").Value = 124912323
$ws.Range("G38").Value = 0.742

$ws.Range("B39").Value = "ARB"
$ws.Range("C39").Value = "Arbitrum"
$ws.Range("D39").Value = 1.12
$ws.Range("E39").Value = 3243118839
$ws.Range("F39").Value = 148252828
$ws.Range("G39").Value = -0.29229

$ws.Range("B40").Value = "FIL"
$ws.Range("C40").Value = "Filecoin"
$ws.Range("D40").Value = 5.78
$ws.Range("E40").Value = 3224187823
$ws.Range("F40").Value = 94334421
$ws.Range("G40").Value = 0.32729

$ws.Range("B41").Value = "MNT"
$ws.Range("C41").Value = "Mantle"
$ws.Range("D41").Value = 0.986084
$ws.Range("E41").Value = 3216433285
$ws.Range("F41").Value = 47240688
$ws.Range("G41").Value = -0.831

$ws.Range("B42").Value = "XLM"
$ws.Range("C42").Value = "Stellar"
$ws.Range("D42").Value = 0.105776
$ws.Range("E42").Value = 3069275231
$ws.Range("F42").Value = 31852429
$ws.Range("G42").Value = -0.73455

$ws.Range("B43").Value = "USDE"
$ws.Range("C43").Value = "Ethena USDe"
$ws.Range("D43").Value = 1.001
$ws.Range("E43").Value = 3017705975
$ws.Range("F43").Value = 34794630
$ws.Range("G43").Value = -0.00293

$ws.Range("B44").Value = "CRO"
$ws.Range("C44").Value = "Cronos"
$ws.Range("D44").Value = 0.112735
$ws.Range("E44").Value = 3016490752
$ws.Range("F44").Value = 6737163
$ws.Range("G44").Value = -0.34543

$ws.Range("B45").Value = "FDUSD"
$ws.Range("C45").Value = "First Digital USD"
$ws.Range("D45").Value = 0.9978320000000001
$ws.Range("E45").Value = 2899835928
$ws.Range("F45").Value = 2578352222
$ws.Range("G45").Value = -0.07239

$ws.Range("B46").Value = "GRT"
$ws.Range("C46").Value = "The Graph"
$ws.Range("D46").Value = 0.298861
$ws.Range("E46").Value = 2841885664
$ws.Range("F46").Value = 45098839
$ws.Range("G46").Value = 0.78516

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "OKB"
$ws.Range("D47").Value = 47.41
$ws.Range("E47").Value = 2840618927
$ws.Range("F47").Value = 3516527
$ws.Range("G47").Value = -1.45975

$ws.Range("B48").Value = "AR"
$ws.Range("C48").Value = "Arweave"
$ws.Range("D48").Value = 43.17
$ws.Range("E48").Value = 2833733478
$ws.Range("F48").Value = 92988423
$ws.Range("G48").Value = -0.892

$ws.Range("B49").Value = "XMR"
$ws.Range("C49").Value = "Monero"
$ws.Range("D49").Value = 151.03
$ws.Range("E49").Value = 2741159306
$ws.Range("F49").Value = 47465868
$ws.Range("G49").Value = 0.21232

$ws.Range("B50").Value = "TAO"
$ws.Range("C50").Value = "Bittensor"
$ws.Range("D50").Value = 394.84
$ws.Range("E50").Value = 2727113731
$ws.Range("F50").Value = 28517117
$ws.Range("G50").Value = 1.34382

$ws.Range("B51").Value = "STX"
$ws.Range("C51").Value = "Stacks"
$ws.Range("D51").Value = 1.86
$ws.Range("E51").Value = 2722040538
$ws.Range("F51").Value = 21398338
$ws.Range("G51").Value = 2.65289

Write-Output "done"